$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("M2").Value = 1.011219333333333
$ws.Range("N2").Value = 3.033658
$ws.Range("O2").Value = 0.1971977600385235
$ws.Range("P2").Value = 0.1971977600385236
$ws.Range("Q2").Value = 0.1171854895164444
$ws.Range("R2").Value = 1.054669405648
$ws.Range("S2").Value = 0.1971977600385235
$ws.Range("T2").Value = 0.1971977600385236

# --- Update existing row 3 values ---
$ws.Range("M3").Value = 3.237665999999999
$ws.Range("N3").Value = 9.712997999999999
$ws.Range("O3").Value = 0.6313768555515021
$ws.Range("P3").Value = 0.6313768555515022
$ws.Range("Q3").Value = 0.375198003632
$ws.Range("S3").Value = 0.6313768555515021
$ws.Range("T3").Value = 0.6313768555515022

# --- Update existing row 4 values ---
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7485606666666667
$ws.Range("N4").Value = 2.245682
$ws.Range("O4").Value = 0.1459767251808977
$ws.Range("P4").Value = 0.1459767251808977
$ws.Range("Q4").Value = 0.08674720237688889
$ws.Range("R4").Value = 0.7807248213920001
$ws.Range("S4").Value = 0.1459767251808977
$ws.Range("T4").Value = 0.1459767251808977

# --- Add new row 5 ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "a"
$ws.Range("C5").Value = "Mc5r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1158853333333333
$ws.Range("H5").Value = 0.347656
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1304993333333333
$ws.Range("N5").Value = 0.391498
$ws.Range("O5").Value = 0.02544865922907654
$ws.Range("P5").Value = 0.02544865922907655
$ws.Range("Q5").Value = 0.01512295874311111
$ws.Range("R5").Value = 0.136106628688
$ws.Range("S5").Value = 0.02544865922907654
$ws.Range("T5").Value = 0.02544865922907655

$wb.Save()
